$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    while ($r.Find.Execute($find, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)) {
        $r.Text = $replace
        $r.Collapse(0)
    }
}

Replace-Text "FELICITAS M. SUMAGUI" "TERESITA P. RIÑO"
Replace-Text "City Social Welfare Development Office" "Vice Mayor's Office"
Replace-Text "120.916" "115.000"
Replace-Text "68,843.35" "65,475.08"
